$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-03-20 Thursday" "2025-03-21 Friday"

Replace-Text "544×8=" "996×9="
Replace-Text "863×6=" "150×7="
Replace-Text "438×5=" "436×8="
Replace-Text "768×6=" "209×5="
Replace-Text "330×8=" "390×8="

Replace-Text "512×2=" "237×6="
Replace-Text "754×9=" "637×8="
Replace-Text "870×6=" "303×4="
Replace-Text "205×4=" "229×7="
Replace-Text "281×3=" "863×7="

Replace-Text "825×8=" "848×6="
Replace-Text "224×2=" "306×2="
Replace-Text "595×5=" "630×9="
Replace-Text "302×6=" "282×3="
Replace-Text "982×6=" "305×7="

Replace-Text "332×9=" "966×8="
Replace-Text "257×5=" "561×5="
Replace-Text "352×7=" "550×5="
Replace-Text "514×3=" "779×8="
Replace-Text "878×2=" "159×4="

Replace-Text "905×6=" "275×6="
Replace-Text "591×9=" "880×8="
Replace-Text "932×5=" "183×3="
Replace-Text "285×7=" "879×4="
Replace-Text "340×4=" "737×4="
